# Insert a new data row at row 80 (pushes the existing rows 80..138 down to 81..139)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with the new "Acelga" record.
$ws.Range("A80").Value = 4
$ws.Range("B80").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C80").Value = "Los Lagos"
$ws.Range("D80").Value = 44574
$ws.Range("E80").Value = 10
$ws.Range("F80").Value = 100112009
$ws.Range("G80").Value = "Acelga"
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 40
$ws.Range("K80").Value = 10000
$ws.Range("L80").Value = 10000
$ws.Range("M80").Value = 10000
$ws.Range("N80").Value = "$/docena de atados (12 kilos)"
$ws.Range("O80").Value = "Región de La Araucanía"
$ws.Range("P80").Value = 833
$ws.Range("Q80").Value = 12
$ws.Range("R80").Value = "Hortaliza"
